$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E text-valued cells remain stored as text (not auto-converted to numbers)

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '24.788.30'

$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.76%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.702.46'

$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.36%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.32%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '317.31'

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.07%  '

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.33%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3956'

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.33%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4086'

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -1.34%  '

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.35%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '52.66'

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.04%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08917'

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.74%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.709'

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +6.79%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '24.38'

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +4.74%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.157'

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.19%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001333'

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.36%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.707.92'

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.75%  '

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.15%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07129'

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.81%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '20.06'

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.82%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.227'

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +3.95%  '

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.69%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.61'

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +3.04%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '24.791.37'

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.78%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.107'

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.71%  '

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.09%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '23.01'

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.66%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.359'

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +25.28%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '165.17'

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.51%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '139.70'

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +2.45%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.193'

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.105'

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +13.46%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09170'

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +6.60%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.083'

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.54%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.03052'

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +11.41%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.2816'

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +2.80%  '

$ws.Range('B37').Value = 'FraxShare'

$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '11.09'

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -4.21%  '

$ws.Range('B38').Value = 'WEMIXTOKEN'

$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.966'

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +1.91%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '14.58'

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.61%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.09314'

$ws.Range('B41').Value = 'TrustWalletToken'

$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.483'

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.06%  '

$ws.Range('B42').Value = 'TheSandbox'

$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.7832'

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +2.14%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '16.30'

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +4.33%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.637'

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +3.41%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.7271'

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.32%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.251'

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.73%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.362'

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +3.12%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.003'

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.29%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '141.25'

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +0.22%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '93.39'

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +4.75%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.08067'

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.00%  '
